$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 588.7778
$ws.Range("I52").Value = 133.16667
$ws.Range("J52").Value = 1500
$ws.Range("K52").Value = 399.50001
$ws.Range("L52").Value = 4500
$ws.Range("M52").Value = -239.50001
$ws.Range("N52").Value = -4820

$ws.Range("H74").Value = 5716.6665
$ws.Range("I74").Value = 5233.3335
$ws.Range("J74").Value = 6200
$ws.Range("K74").Value = 5233.3335
$ws.Range("L74").Value = 6200
$ws.Range("M74").Value = -4297.3335
$ws.Range("N74").Value = -8072

$ws.Range("H77").Value = 5716.6665
$ws.Range("I77").Value = 5233.3335
$ws.Range("J77").Value = 6200
$ws.Range("K77").Value = 26166.6675
$ws.Range("L77").Value = 31000
$ws.Range("M77").Value = -21486.6675
$ws.Range("N77").Value = -40360

$ws.Range("H111").Value = 11558.786
$ws.Range("I111").Value = 25917.4
$ws.Range("J111").Value = 3581.7778
$ws.Range("K111").Value = 77752.20000000001
$ws.Range("L111").Value = 10745.3334
$ws.Range("M111").Value = -74685.20000000001
$ws.Range("N111").Value = -16879.3334

$ws.Range("H123").Value = 20238
$ws.Range("J123").Value = 20238
$ws.Range("L123").Value = 20238
$ws.Range("N123").Value = -30038

$ws.Range("H129").Value = 833.0492
$ws.Range("J129").Value = 878.2963
$ws.Range("L129").Value = 2634.8889
$ws.Range("N129").Value = -12634.8889

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 32509.281
$ws.Range("I2").Value = 1193.6842
$ws.Range("J2").Value = 78278.234
$ws.Range("K2").Value = 1193.6842
$ws.Range("L2").Value = 78278.234
$ws.Range("M2").Value = -1080.6842
$ws.Range("N2").Value = -78504.234

$ws.Range("H32").Value = 28199.39
$ws.Range("I32").Value = 4636.755
$ws.Range("J32").Value = 141728.45
$ws.Range("K32").Value = 4636.755
$ws.Range("L32").Value = 141728.45
$ws.Range("M32").Value = -4349.755
$ws.Range("N32").Value = -142302.45

$ws.Range("H61").Value = 2133.0967
$ws.Range("I61").Value = 1078.1538
$ws.Range("J61").Value = 2895
$ws.Range("K61").Value = 1078.1538
$ws.Range("L61").Value = 2895
$ws.Range("M61").Value = -866.1538
$ws.Range("N61").Value = -3319

$ws.Range("H116").Value = 32509.281
$ws.Range("I116").Value = 1193.6842
$ws.Range("J116").Value = 78278.234
$ws.Range("K116").Value = 1193.6842
$ws.Range("L116").Value = 78278.234
$ws.Range("M116").Value = 1100.3158
$ws.Range("N116").Value = -82866.234

$ws.Range("H136").Value = 2133.0967
$ws.Range("I136").Value = 1078.1538
$ws.Range("J136").Value = 2895
$ws.Range("K136").Value = 3234.4614
$ws.Range("L136").Value = 8685
$ws.Range("M136").Value = -684.4614000000001
$ws.Range("N136").Value = -13785

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 32509.281
$ws.Range("I3").Value = 1193.6842
$ws.Range("J3").Value = 78278.234
$ws.Range("K3").Value = 1193.6842
$ws.Range("L3").Value = 78278.234
$ws.Range("M3").Value = -1079.6842
$ws.Range("N3").Value = -78506.234

$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("M141").Value = ""

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 4566
$ws.Range("I32").Value = 0
$ws.Range("K32").Value = 0
$ws.Range("M32").Value = ""

$ws.Range("H34").Value = 2227.0908
$ws.Range("J34").Value = 2688.6667
$ws.Range("L34").Value = 8066.000100000001
$ws.Range("N34").Value = -8234.000100000001

$ws.Range("H109").Value = 3115.652
$ws.Range("I109").Value = 458
$ws.Range("J109").Value = 4278.375
$ws.Range("K109").Value = 1374
$ws.Range("L109").Value = 12835.125
$ws.Range("M109").Value = -334
$ws.Range("N109").Value = -14915.125

$ws.Range("H131").Value = 817.92
$ws.Range("I131").Value = 432
$ws.Range("J131").Value = 838.23157
$ws.Range("K131").Value = 1296
$ws.Range("L131").Value = 2514.69471
$ws.Range("M131").Value = 3744
$ws.Range("N131").Value = -12594.69471

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 58148.184
$ws.Range("I70").Value = 82151.38
$ws.Range("J70").Value = 6141.25
$ws.Range("K70").Value = 82151.38
$ws.Range("L70").Value = 6141.25
$ws.Range("M70").Value = -81881.38
$ws.Range("N70").Value = -6681.25

$ws.Range("H73").Value = 58148.184
$ws.Range("I73").Value = 82151.38
$ws.Range("J73").Value = 6141.25
$ws.Range("K73").Value = 82151.38
$ws.Range("L73").Value = 6141.25
$ws.Range("M73").Value = -81215.38
$ws.Range("N73").Value = -8013.25

$ws.Range("H135").Value = 31472.8
$ws.Range("J135").Value = 31472.8
$ws.Range("L135").Value = 31472.8
$ws.Range("N135").Value = -41612.8

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").Value = ""

$ws.Range("H138").Value = 59947.5
$ws.Range("J138").Value = 65395
$ws.Range("L138").Value = 65395
$ws.Range("N138").Value = -75675

$ws.Range("H139").Value = 44600
$ws.Range("J139").Value = 44600
$ws.Range("L139").Value = 44600
$ws.Range("N139").Value = -54880

$ws.Range("H140").Value = 70000
$ws.Range("J140").Value = 70000
$ws.Range("L140").Value = 70000
$ws.Range("N140").Value = -80360

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 39330
$ws.Range("J16").Value = 39330
$ws.Range("L16").Value = 39330
$ws.Range("N16").Value = -39914

$ws.Range("H56").Value = 32585.9
$ws.Range("J56").Value = 35841.555
$ws.Range("L56").Value = 35841.555
$ws.Range("N56").Value = -37269.555

$ws.Range("H132").Value = 7220.1113
$ws.Range("I132").Value = 15668
$ws.Range("J132").Value = 2996.1667
$ws.Range("K132").Value = 47004
$ws.Range("L132").Value = 8988.500100000001
$ws.Range("M132").Value = -44474
$ws.Range("N132").Value = -14048.5001

$ws.Range("H135").Value = 49900
$ws.Range("J135").Value = 49900
$ws.Range("L135").Value = 49900
$ws.Range("N135").Value = -60040

$ws.Range("H136").Value = 1574.4706
$ws.Range("I136").Value = 591.3333
$ws.Range("J136").Value = 1928.4
$ws.Range("K136").Value = 1773.9999
$ws.Range("L136").Value = 5785.200000000001
$ws.Range("M136").Value = 776.0001
$ws.Range("N136").Value = -10885.2

$ws.Range("H137").Value = 49900
$ws.Range("J137").Value = 49900
$ws.Range("L137").Value = 49900
$ws.Range("N137").Value = -60100

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = ""

$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").Value = ""

$ws.Range("H140").Value = 56997.5
$ws.Range("J140").Value = 56997.5
$ws.Range("L140").Value = 56997.5
$ws.Range("N140").Value = -67357.5

$ws.Range("H141").Value = 65715
$ws.Range("J141").Value = 65715
$ws.Range("L141").Value = 65715
$ws.Range("N141").Value = -76075
